# Regenerate s_val data to filter save games.
# Updates columns B (TB), C (d2S), D (K), E (IP), G (sum) for rows 2-14.
# Column A (date) and F (Win) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    3  = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 4.429675500412797 }
    4  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 }
    5  = @{ B = 0.6753301551942219; C = 0.3127903958511391; D = 0.8054896365839992; E = 8.660232485948974;  G = 10.45384267357833 }
    6  = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 3.781711156805759 }
    7  = @{ B = 0.127881588408715;  C = 0.3127903958511391; D = 0.8054896365839992; E = 8.660232485948974;  G = 9.906394106792828 }
    8  = @{ B = 0.01514828764759746; C = 0.04240448674262143; D = 0.1575252929769615; E = 0.496779210170732; G = 0.7118572775379124 }
    9  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    10 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    11 = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 2.997429241610044 }
    12 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 4.429675500412797 }
    13 = @{ B = 0.6753301551942219; C = 0.002777888934908601; D = 0.8054896365839992; E = 0.496779210170732; G = 1.980376890883862 }
    14 = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732;  G = 6.740334628841572 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
